$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the content of three pairs of rows in the JobPlanning
# sheet (rows 2<->3, 5<->7, 6<->8), columns A:M. Using Copy +
# PasteSpecial(xlPasteValues) instead of a plain .Value assignment keeps
# each destination cell's original style AND its original cell type
# (numeric-looking strings like "103"/"124" stay text/shared-strings
# instead of being re-interpreted as numbers).

function Swap-Rows($rowA, $rowB) {
    $rangeA = "A" + $rowA + ":M" + $rowA
    $rangeB = "A" + $rowB + ":M" + $rowB
    $scratch = "A200:M200"

    $ws.Range($rangeA).Copy()
    $ws.Range($scratch).PasteSpecial(-4163)

    $ws.Range($rangeB).Copy()
    $ws.Range($rangeA).PasteSpecial(-4163)

    $ws.Range($scratch).Copy()
    $ws.Range($rangeB).PasteSpecial(-4163)

    $ws.Range($scratch).ClearContents()
}

Swap-Rows 2 3
Swap-Rows 5 7
Swap-Rows 6 8
